$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns retain their original text formatting so that
# numeric-looking strings (e.g. "303.50", "6.130") are not auto-converted to numbers.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.360.26"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.23%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.625.64"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.01%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.59%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.003"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.59%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "303.50"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.29%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3765"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.51%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "51.98"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.99%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3618"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.20%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.232"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.35%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08079"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.32%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.59%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.65"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.24%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.553"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.54%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001246"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.77%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.215"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.92%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.625.25"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.23%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "93.38"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.56%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06922"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.34%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.94"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.25%  "

# Row 21
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.31%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.454"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.90%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "23.360.30"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.23%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.73"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.94%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.219"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.19%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.416"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.04%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.09"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.31%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "148.97"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.33%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.301"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.34%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.55"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.48%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.301"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -4.54%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.805.79"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.47%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.776"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.49%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "10.92"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +4.17%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9484"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.25%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02835"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.46%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2529"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.11%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.130"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.68%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.08811"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.35%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.07114"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -4.29%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.361"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -3.06%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7042"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.75%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.11"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.65%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.33"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.19%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6464"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.34%  "

# Row 46
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.317"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.66%  "

# Row 47
$ws.Range("B47").Value = "Frax"
$ws.Range("C47").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.001"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.51%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.982"

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07959"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.51%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.206"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.38%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "126.09"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -4.06%  "
